# Fix response-template sample data on sheet "Лист1":
# the sample password in A3 was wrong ("0u0p4m4U") - update it to the
# correct template word ("0u0p4M4u@!") and, since that value now reads
# like an email-ish token (it contains "@"), turn it into the hyperlink
# that Excel creates for such tokens (this also applies the built-in
# "Hyperlink" cell style - underlined, theme color 10 - to A3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Correct the sample template word in A3.
$ws.Range("A3").Value = "0u0p4M4u@!"

# 2) Turn it into a hyperlink (mailto:, matching the text itself) and let
#    Excel apply its built-in Hyperlink style to the cell.
$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:0u0p4M4u@!")

# 3) Leave the selection on A3, like in the saved workbook.
[void]$ws.Range("A3").Select()
